$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text in the source data -
# some look like ordinary decimals (e.g. "1.016") that Excel would otherwise
# silently reinterpret as numbers, and others use "." as a thousands
# separator (e.g. "27.482.62") which is not a legal number at all. Force
# every cell whose price we are about to rewrite to Text format first so the
# new value is stored the same way the original was (as a string).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.482.62"
$ws.Range("E2").Value = "  +2.27%  "
$ws.Range("D3").Value = "1.872.85"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("D4").Value = "1.016"
$ws.Range("E4").Value = "  +0.91%  "
$ws.Range("D5").Value = "313.33"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").Value = "1.013"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("D7").Value = "0.4782"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").Value = "0.3770"
$ws.Range("E8").Value = "  +2.68%  "
$ws.Range("D9").Value = "0.07377"
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("D10").Value = "0.9367"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").Value = "20.69"
$ws.Range("E11").Value = "  +5.56%  "
$ws.Range("D12").Value = "0.07842"
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("D13").Value = "1.873.39"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "5.446"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").Value = "6.588"
$ws.Range("E15").Value = "  +2.92%  "
$ws.Range("D16").Value = "90.91"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "1.017"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").Value = "0.000008929"
$ws.Range("E18").Value = "  +3.37%  "
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").Value = "14.90"
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").Value = "27.529.48"
$ws.Range("E21").Value = "  +2.35%  "
$ws.Range("D22").Value = "5.133"
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("D23").Value = "10.73"
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("D24").Value = "1.966"
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("D25").Value = "153.91"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("D27").Value = "2.019"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("D28").Value = "115.92"
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("D29").Value = "4.996"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").Value = "0.08933"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").Value = "3.344"
$ws.Range("E31").Value = "  +1.81%  "
$ws.Range("D32").Value = "1.217"
$ws.Range("E32").Value = "  +3.76%  "
$ws.Range("D33").Value = "4.619"
$ws.Range("E33").Value = "  +2.98%  "
$ws.Range("D34").Value = "0.7517"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("D35").Value = "2.688"
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("D36").Value = "0.02067"
$ws.Range("E36").Value = "  +6.32%  "
$ws.Range("D37").Value = "1.118"
$ws.Range("E37").Value = "  +2.70%  "
$ws.Range("D38").Value = "0.05302"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").Value = "3.006"
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("D40").Value = "0.5348"
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("D41").Value = "7.089"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").Value = "8.426"
$ws.Range("E43").Value = "  +2.58%  "
$ws.Range("D44").Value = "10.61"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").Value = "0.4833"
$ws.Range("E45").Value = "  +2.38%  "
$ws.Range("D46").Value = "1.014"
$ws.Range("E46").Value = "  +0.79%  "
$ws.Range("D47").Value = "1.662"
$ws.Range("E47").Value = "  +3.67%  "
$ws.Range("D48").Value = "103.05"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("D49").Value = "67.31"
$ws.Range("E49").Value = "  +2.94%  "
$ws.Range("D50").Value = "0.06094"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("D51").Value = "0.8970"
$ws.Range("E51").Value = "  +1.39%  "
